$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.159.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.846.61'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3705'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07377'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8839'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07932'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.91'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.871.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.373'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.593'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008949'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.15%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  +3.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.187.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.142'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.127.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.867'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.069'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.140'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08887'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.973'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7410'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.469'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.142'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.547'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.080'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05272'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01951'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.969'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.092'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5175'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1637'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.273'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4857'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.005'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.82'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.631'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06231'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.10%  '
